# Weekly update: insert a new "Fruta, Terminal Hortofrutícola Agro Chillán -
# Manzana" price record at row 603, pushing the existing rows 603:636 down to
# 604:637 (dimension grows from A1:T636 to A1:T637).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 603; Excel shifts rows 603:636
# down to 604:637 and carries the row-603 formatting (incl. the date style
# on column D) down with them / onto the new row.
$ws.Rows.Item(603).Insert()

# Populate the newly inserted row 603 with this week's record.
$ws.Range("A603").Value = 7
$ws.Range("B603").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C603").Value = "Ñuble"
$ws.Range("D603").Value = 44610
$ws.Range("E603").Value = 16
$ws.Range("F603").Value = "Fruta"
$ws.Range("G603").Value = 100104
$ws.Range("H603").Value = "Frutos de pepita"
$ws.Range("I603").Value = 100104002
$ws.Range("J603").Value = "Manzana"
$ws.Range("K603").Value = "Royal Gala"
$ws.Range("L603").Value = "Segunda"
$ws.Range("M603").Value = 60
$ws.Range("N603").Value = 9000
$ws.Range("O603").Value = 9000
$ws.Range("P603").Value = 9000
$ws.Range("Q603").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R603").Value = "Región de O'Higgins"
$ws.Range("S603").Value = 562
$ws.Range("T603").Value = 16
